$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder / card number
$ws.Range("C2").Value = "Hartmut"
# Force the 16-digit card number to stay text (matches source which stores it
# as an inline string), then restore the original cell formatting.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 30.09.2023"

# Row 6
$ws.Range("B6").Value = "04.10."
$ws.Range("C6").Value = "05.10."
$ws.Range("D6").Value = "MCDONALDS Geithain"
$ws.Range("E6").Value = "30,22-"

# Row 7
$ws.Range("B7").Value = "07.10."
$ws.Range("C7").Value = "08.10."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU KLJTJI"
$ws.Range("E7").Value = "239,99-"

# Row 8
$ws.Range("B8").Value = "11.10."
$ws.Range("C8").Value = "12.10."
$ws.Range("D8").Value = "KARTENZ./11.10 LIDL RO"
$ws.Range("E8").Value = "55,37-"

# Row 9 (previously empty placeholder row) - pull in matching formats from row 8
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("B9").Value = "12.10."
$ws.Range("C9").Value = "13.10."
$ws.Range("D9").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Value = "68,78-"

# Row 10 (previously empty placeholder row)
$ws.Range("B8:E8").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)
$ws.Range("B10").Value = "16.10."
$ws.Range("C10").Value = "17.10."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-62320686"
$ws.Range("E10").Value = "57,23-"

# Row 11 (previously empty placeholder row)
$ws.Range("B8:E8").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)
$ws.Range("B11").Value = "18.10."
$ws.Range("C11").Value = "19.10."
$ws.Range("D11").Value = "PAYPAL DTORLN"
$ws.Range("E11").Value = "51,97-"

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 21.10.2023"
$ws.Range("E12").Value = "503,56-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 27.10.2023"
